$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test run flipped the pass/fail result for the second "tutorial"
# login attempt: C2 (was FAIL) is now PASS, and C3 (was PASS) is now FAIL.
# Swap the two cells' full state (value + fill) so the PASS/FAIL text
# moves together with its red/green status colour.
$c2 = $ws.Range("C2")
$c3 = $ws.Range("C3")

# Use a scratch cell well outside the used range as a holding spot so we
# can do a clean 3-way swap (C2 -> tmp -> C2 becomes C3 -> C3 becomes tmp).
$tmp = $ws.Range("A100")

$c2.Cut($tmp)
$c3.Cut($c2)
$tmp.Cut($c3)

# Remove the scratch cell again so it doesn't leave a stray entry / grow
# the sheet's used range.
$tmp.Delete(-4162)  # xlUp
